$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 224
$ws.Range("I2").Value = 560
$ws.Range("J2").Value = 2433
$ws.Range("K2").Value = 16
$ws.Range("L2").Value = 677
$ws.Range("M2").Value = 32
$ws.Range("N2").Value = 409
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 9
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 31
$ws.Range("S2").Value = 264
$ws.Range("T2").Value = 439
$ws.Range("U2").Value = 30
$ws.Range("V2").Value = 3671
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 3837
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 66
$ws.Range("AA2").Value = 27
